$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 9695
$ws.Range("J16").Value = 9695
$ws.Range("L16").Value = 9695
$ws.Range("N16").Value = -10155

$ws.Range("H28").Value = 450.3889
$ws.Range("I28").Value = 464.7857
$ws.Range("K28").Value = 464.7857
$ws.Range("M28").Value = 20.21429999999998

$ws.Range("H33").Value = 4847094
$ws.Range("I33").Value = 6230374
$ws.Range("J33").Value = 5614.75
$ws.Range("K33").Value = 6230374
$ws.Range("L33").Value = 5614.75
$ws.Range("M33").Value = -6230145
$ws.Range("N33").Value = -6072.75

$ws.Range("H41").Value = 1052.7646
$ws.Range("I41").Value = 520
$ws.Range("K41").Value = 520
$ws.Range("M41").Value = -80

$ws.Range("H70").Value = 2776
$ws.Range("I70").Value = 857.6667
$ws.Range("K70").Value = 2573.0001
$ws.Range("M70").Value = -2303.0001

$ws.Range("H73").Value = 2776
$ws.Range("I73").Value = 857.6667
$ws.Range("K73").Value = 2573.0001
$ws.Range("M73").Value = -1637.0001

$ws.Range("H86").Value = 666666700
$ws.Range("I86").Value = 1000000000
$ws.Range("J86").Value = 333333340
$ws.Range("K86").Value = 1000000000
$ws.Range("L86").Value = 333333340
$ws.Range("M86").Value = -999998877
$ws.Range("N86").Value = -333335586

$ws.Range("H89").Value = 666666700
$ws.Range("I89").Value = 1000000000
$ws.Range("J89").Value = 333333340
$ws.Range("K89").Value = 5000000000
$ws.Range("L89").Value = 1666666700
$ws.Range("M89").Value = -4999994384
$ws.Range("N89").Value = -1666677932

$ws.Range("H132").Value = 1630.32
$ws.Range("I132").Value = 1677.3158
$ws.Range("K132").Value = 5031.9474
$ws.Range("M132").Value = -2501.9474

$ws.Range("H135").Value = 350.55554
$ws.Range("J135").Value = 197.5
$ws.Range("L135").Value = 1777.5
$ws.Range("N135").Value = -6847.5

$ws.Range("H137").Value = 3712587.8
$ws.Range("I137").Value = 5993.95
$ws.Range("J137").Value = 14302856
$ws.Range("K137").Value = 17981.85
$ws.Range("L137").Value = 42908568
$ws.Range("M137").Value = -15431.85
$ws.Range("N137").Value = -42913668

$ws.Range("H138").Value = 5689.3486
$ws.Range("J138").Value = 2425.0322
$ws.Range("L138").Value = 7275.096600000001
$ws.Range("N138").Value = -17555.0966

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 824.1667
$ws.Range("I4").Value = 849.2
$ws.Range("K4").Value = 849.2
$ws.Range("M4").Value = -733.2

$ws.Range("H32").Value = 306506
$ws.Range("I32").Value = 400627.56
$ws.Range("K32").Value = 400627.56
$ws.Range("M32").Value = -400340.56

$ws.Range("H88").Value = 2107
$ws.Range("I88").Value = 1449.5
$ws.Range("K88").Value = 1449.5
$ws.Range("M88").Value = -1043.5

$ws.Range("H91").Value = 2107
$ws.Range("I91").Value = 1449.5
$ws.Range("K91").Value = 1449.5
$ws.Range("M91").Value = -45.5

$ws.Range("H97").Value = 9309.200000000001
$ws.Range("I97").Value = 10688.546
$ws.Range("K97").Value = 10688.546
$ws.Range("M97").Value = -10192.546

$ws.Range("H102").Value = 3003.3333
$ws.Range("I102").Value = 3003.3333
$ws.Range("K102").Value = 3003.3333
$ws.Range("M102").Value = -1381.3333

$ws.Range("H132").Value = 3535.3215
$ws.Range("I132").Value = 2592.8
$ws.Range("J132").Value = 4058.9443
$ws.Range("K132").Value = 7778.400000000001
$ws.Range("L132").Value = 12176.8329
$ws.Range("M132").Value = -5248.400000000001
$ws.Range("N132").Value = -17236.8329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5036.5
$ws.Range("I86").Value = 1457.2
$ws.Range("J86").Value = 11002
$ws.Range("K86").Value = 1457.2
$ws.Range("L86").Value = 11002
$ws.Range("M86").Value = -334.2
$ws.Range("N86").Value = -13248

$ws.Range("H89").Value = 5036.5
$ws.Range("I89").Value = 1457.2
$ws.Range("J89").Value = 11002
$ws.Range("K89").Value = 7286
$ws.Range("L89").Value = 55010
$ws.Range("M89").Value = -1670
$ws.Range("N89").Value = -66242

$ws.Range("H107").Value = 40995.4
$ws.Range("I107").Value = 50494.25
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 50494.25
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -48574.25
$ws.Range("N107").Value = -6840

$ws.Range("H134").Value = 24326214
$ws.Range("I134").Value = 1877.2106
$ws.Range("K134").Value = 5631.6318
$ws.Range("M134").Value = -3096.6318

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 204.53847
$ws.Range("I7").Value = 303.625
$ws.Range("K7").Value = 303.625
$ws.Range("M7").Value = -190.625

$ws.Range("H22").Value = 1544.6666
$ws.Range("I22").Value = 3001
$ws.Range("J22").Value = 816.5
$ws.Range("K22").Value = 3001
$ws.Range("L22").Value = 816.5
$ws.Range("M22").Value = -2651
$ws.Range("N22").Value = -1516.5

$ws.Range("H107").Value = 2162.1365
$ws.Range("I107").Value = 2210.1052
$ws.Range("J107").Value = 1858.3334
$ws.Range("K107").Value = 2210.1052
$ws.Range("L107").Value = 1858.3334
$ws.Range("M107").Value = -290.1052
$ws.Range("N107").Value = -5698.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 743.4761999999999
$ws.Range("J12").Value = 1043.0714
$ws.Range("L12").Value = 3129.2142
$ws.Range("N12").Value = -3475.2142

$ws.Range("H113").Value = 514.8
$ws.Range("J113").Value = 514.8
$ws.Range("L113").Value = 1544.4
$ws.Range("N113").Value = -5884.4

$ws.Range("H124").Value = 3347.9
$ws.Range("I124").Value = 2652.6667
$ws.Range("K124").Value = 7958.000100000001
$ws.Range("M124").Value = -3048.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 38475930
$ws.Range("I80").Value = 14848.75
$ws.Range("J80").Value = 55569740
$ws.Range("K80").Value = 14848.75
$ws.Range("L80").Value = 55569740
$ws.Range("M80").Value = -13850.75
$ws.Range("N80").Value = -55571736

$ws.Range("H83").Value = 38475930
$ws.Range("I83").Value = 14848.75
$ws.Range("J83").Value = 55569740
$ws.Range("K83").Value = 74243.75
$ws.Range("L83").Value = 277848700
$ws.Range("M83").Value = -69251.75
$ws.Range("N83").Value = -277858684

$ws.Range("H97").Value = 1026.4546
$ws.Range("I97").Value = 479.2
$ws.Range("K97").Value = 479.2
$ws.Range("M97").Value = 16.80000000000001

$ws.Range("H107").Value = 125999.375
$ws.Range("J107").Value = 1199.5
$ws.Range("L107").Value = 1199.5
$ws.Range("N107").Value = -5039.5

$ws.Range("H113").Value = 9997
$ws.Range("J113").Value = 9997
$ws.Range("L113").Value = 9997
$ws.Range("N113").Value = -14337

$ws.Range("H132").Value = 886627.9
$ws.Range("I132").Value = 1076.0588
$ws.Range("K132").Value = 3228.1764
$ws.Range("M132").Value = -698.1764000000003

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1636284.4
$ws.Range("J40").Value = 2014.4
$ws.Range("L40").Value = 2014.4
$ws.Range("N40").Value = -2286.4

$ws.Range("H55").Value = 933.2121
$ws.Range("I55").Value = 1136.7368
$ws.Range("J55").Value = 657
$ws.Range("K55").Value = 1136.7368
$ws.Range("L55").Value = 657
$ws.Range("M55").Value = -963.7367999999999
$ws.Range("N55").Value = -1003

$ws.Range("H61").Value = 3562.4666
$ws.Range("I61").Value = 2493.5557
$ws.Range("J61").Value = 5165.8335
$ws.Range("K61").Value = 2493.5557
$ws.Range("L61").Value = 5165.8335
$ws.Range("M61").Value = -2291.5557
$ws.Range("N61").Value = -5569.8335

$ws.Range("H113").Value = 3562.4666
$ws.Range("I113").Value = 2493.5557
$ws.Range("J113").Value = 5165.8335
$ws.Range("K113").Value = 2493.5557
$ws.Range("L113").Value = 5165.8335
$ws.Range("M113").Value = -323.5556999999999
$ws.Range("N113").Value = -9505.833500000001

$ws.Range("H132").Value = 4583.6924
$ws.Range("I132").Value = 2519.8
$ws.Range("K132").Value = 7559.400000000001
$ws.Range("M132").Value = -5029.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 92763.37
$ws.Range("I81").Value = 2199.8572
$ws.Range("J81").Value = 251249.5
$ws.Range("K81").Value = 4399.7144
$ws.Range("L81").Value = 502499
$ws.Range("M81").Value = -3338.7144
$ws.Range("N81").Value = -504621

$ws.Range("H84").Value = 92763.37
$ws.Range("I84").Value = 2199.8572
$ws.Range("J84").Value = 251249.5
$ws.Range("K84").Value = 21998.572
$ws.Range("L84").Value = 2512495
$ws.Range("M84").Value = -16694.572
$ws.Range("N84").Value = -2523103

$ws.Range("H107").Value = 1787526
$ws.Range("I107").Value = 1512.7142
$ws.Range("J107").Value = 3176647.5
$ws.Range("K107").Value = 4538.142599999999
$ws.Range("L107").Value = 9529942.5
$ws.Range("M107").Value = -2618.142599999999
$ws.Range("N107").Value = -9533782.5

$ws.Range("H113").Value = 502.37036
$ws.Range("I113").Value = 439.6
$ws.Range("J113").Value = 681.7143
$ws.Range("K113").Value = 1318.8
$ws.Range("L113").Value = 2045.1429
$ws.Range("M113").Value = 851.1999999999998
$ws.Range("N113").Value = -6385.1429

$ws.Range("H126").Value = 20835622
$ws.Range("I126").Value = 41668428
$ws.Range("K126").Value = 125005284
$ws.Range("M126").Value = -125002814

$ws.Range("H132").Value = 2001.8214
$ws.Range("I132").Value = 1683.409
$ws.Range("K132").Value = 5050.227000000001
$ws.Range("M132").Value = -2520.227000000001

$ws.Range("H133").Value = 77910
$ws.Range("J133").Value = 77910
$ws.Range("L133").Value = 77910
$ws.Range("N133").Value = -88030

$ws.Range("H136").Value = 786.6111
$ws.Range("I136").Value = 629.2308
$ws.Range("K136").Value = 1887.6924
$ws.Range("M136").Value = 662.3075999999999

$ws.Range("H139").Value = 86500
$ws.Range("J139").Value = 86500
$ws.Range("L139").Value = 86500
$ws.Range("N139").Value = -96780
